$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B holds date-like text values (stored as shared strings, some with a
# quote-prefix so they display verbatim). Update each cell's date text while
# preserving the original cell formatting/quote-prefix behaviour:
#  - rows whose cell kept its quote-prefixed style get a leading apostrophe
#    so Excel keeps treating them as literal text with the prefix flag
#  - row 15 never had the quote-prefix flag, so it is written without one

$ws.Range("B1").Value  = "'19/07/2023"
$ws.Range("B2").Value  = "'19/07/2023"
$ws.Range("B3").Value  = "'19/07/2023"
$ws.Range("B4").Value  = "'19/07/2023"
$ws.Range("B5").Value  = "'19/07/2023"
$ws.Range("B6").Value  = "'19/07/2023"

$ws.Range("B7").Value  = "'17/06/2023"
$ws.Range("B8").Value  = "'17/06/2023"
$ws.Range("B9").Value  = "'17/06/2023"
$ws.Range("B10").Value = "'17/06/2023"
$ws.Range("B11").Value = "'17/06/2023"

$ws.Range("B12").Value = "'2023/06/17"
$ws.Range("B13").Value = "'18/06/2023"
$ws.Range("B14").Value = "'06/17/2023"
$ws.Range("B15").Value = "17/06/2023"
$ws.Range("B16").Value = "'18/06/2023"
$ws.Range("B17").Value = "'18/06/2023"

# Match the new selection left by the edit (single active cell B6)
$ws.Range("B6").Select()
